$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.844.50'
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('D3').Value = '2.082.79'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.32'
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('E6').Value = '  +0.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.23'
$ws.Range('E7').Value = '  +3.63%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.395'
$ws.Range('E9').Value = '  +2.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0789'
$ws.Range('E10').Value = '  +1.78%  '
$ws.Range('E11').Value = '  +1.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.77'
$ws.Range('E12').Value = '  +2.57%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.20'
$ws.Range('E13').Value = '  +0.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.776'
$ws.Range('E14').Value = '  +1.23%  '
$ws.Range('E15').Value = '  +2.81%  '
$ws.Range('D16').Value = '2.069.60'
$ws.Range('E16').Value = '  -1.15%  '
$ws.Range('D17').Value = '37.762.61'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('E18').Value = '  +0.65%  '
$ws.Range('E19').Value = '  +1.27%  '
$ws.Range('D20').Value = '0.0₃0850'
$ws.Range('E20').Value = '  +3.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '228.25'
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('E23').Value = '  -0.50%  '
$ws.Range('E24').Value = '  +1.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.42'
$ws.Range('E25').Value = '  +2.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.17'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.137'
$ws.Range('E27').Value = '  -2.30%  '
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('E29').Value = '  +0.32%  '
$ws.Range('E30').Value = '  +2.03%  '
$ws.Range('E31').Value = '  +2.76%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.75'
$ws.Range('E32').Value = '  +3.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0632'
$ws.Range('E33').Value = '  +1.42%  '
$ws.Range('E34').Value = '  +1.01%  '
$ws.Range('E35').Value = '  +1.65%  '
$ws.Range('E36').Value = '  -0.34%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.41'
$ws.Range('E38').Value = '  +0.49%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0984'
$ws.Range('E39').Value = '  -0.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '99.18'
$ws.Range('E40').Value = '  +1.64%  '
$ws.Range('E41').Value = '  +2.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.89'
$ws.Range('E42').Value = '  -1.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.66'
$ws.Range('E43').Value = '  +6.91%  '
$ws.Range('D44').Value = '1.444.98'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.19'
$ws.Range('E46').Value = '  +3.10%  '
$ws.Range('E47').Value = '  +0.75%  '
$ws.Range('E48').Value = '  +1.25%  '
$ws.Range('E49').Value = '  -0.22%  '
$ws.Range('D50').Value = '2.274.98'
$ws.Range('E50').Value = '  -0.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '46.94'
$ws.Range('E51').Value = '  +1.29%  '
